$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated / new rows for the Efnb1 -> Erbb2 LR-pair table (ECs/FAPs/sCs x ECs/FAPs/sCs).
# Columns: A Sending cluster, B Ligand, C Receptor, D Target cluster, E..T metrics
$rowsData = @(
  @("ECs","Efnb1","Erbb2","ECs",3,1,10.31211433333333,30.936343,0.633340936097251,0.633340936097251,2,0.6666666666666666,2.179771666666667,6.539315,0.2349306639444428,0.2349306639444428,22.47805464722722,202.302491825045,0.1487912066205221,0.1487912066205221),
  @("ECs","Efnb1","Erbb2","FAPs",3,1,10.31211433333333,30.936343,0.633340936097251,0.633340936097251,3,1,4.237801,12.713403,0.4567402255103586,0.4567402255103586,43.70068843391434,393.3061959052291,0.289272281978,0.289272281978),
  @("ECs","Efnb1","Erbb2","sCs",3,1,10.31211433333333,30.936343,0.633340936097251,0.633340936097251,3,1,2.860789,8.582367000000001,0.3083291105451987,0.3083291105451986,29.50078325154234,265.5070492638811,0.1952774474987289,0.1952774474987289),
  @("FAPs","Efnb1","Erbb2","ECs",3,1,4.103438,12.310314,0.2520215719230645,0.2520215719230645,2,0.6666666666666666,2.179771666666667,6.539315,0.2349306639444428,0.2349306639444428,8.944557888323333,80.50102099490999,0.05920759522020769,0.05920759522020767),
  @("FAPs","Efnb1","Erbb2","FAPs",3,1,4.103438,12.310314,0.2520215719230645,0.2520215719230645,3,1,4.237801,12.713403,0.4567402255103586,0.4567402255103586,17.389553659838,156.505982938542,0.1151083895936155,0.1151083895936155),
  @("FAPs","Efnb1","Erbb2","sCs",3,1,4.103438,12.310314,0.2520215719230645,0.2520215719230645,3,1,2.860789,8.582367000000001,0.3083291105451987,0.3083291105451986,11.739070292582,105.651632633238,0.0777055871092413,0.07770558710924126),
  @("sCs","Efnb1","Erbb2","ECs",3,1,1.866538,5.599614,0.1146374919796846,0.1146374919796846,2,0.6666666666666666,2.179771666666667,6.539315,0.2349306639444428,0.2349306639444428,4.068626647156667,36.61763982441,0.02693186210371304,0.02693186210371303),
  @("sCs","Efnb1","Erbb2","FAPs",3,1,1.866538,5.599614,0.1146374919796846,0.1146374919796846,3,1,4.237801,12.713403,0.4567402255103586,0.4567402255103586,7.910016602938001,71.19014942644201,0.05235955393874307,0.05235955393874306),
  @("sCs","Efnb1","Erbb2","sCs",3,1,1.866538,5.599614,0.1146374919796846,0.1146374919796846,3,1,2.860789,8.582367000000001,0.3083291105451987,0.3083291105451986,5.339771378482001,48.05794240633801,0.0353460759372285,0.03534607593722849)
)

$startRow = 2
for ($i = 0; $i -lt $rowsData.Count; $i++) {
    $rowValues = $rowsData[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le $rowValues.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowValues[$c - 1]
    }
}
